# Add a new "FLAG_LETRA" column (G) to the persons sheet, mirroring the
# existing header/data pattern used by columns A-F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "FLAG_LETRA"
$ws.Range("G2").Value = "a"
$ws.Range("G3").Value = "b"
$ws.Range("G4").Value = "c"

# Match column widths of the rest of the header row (auto-sized to content).
$ws.Columns("G:G").AutoFit()

# Leave the selection where Excel would land after typing the last value.
$ws.Range("G5").Select()
